# study2 data updated (Kathy)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the header row so the header moves from row 2 -> row 3
# and the existing data row moves from row 3 -> row 4.
$ws.Rows("2:2").Insert()

# New annotation labels placed in the new row 2 (above the "realism" columns)
$ws.Range("H2").Value = "無"
$ws.Range("J2").Value = "輕"
$ws.Range("L2").Value = "重"
$ws.Range("N2").Value = "輕重"

# New respondent (Kathy) added as row 5, right after the existing respondent (row 4)
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "唐千琳"
$ws.Range("C5").Value = "Kathy"
$ws.Range("D5").Value = "女"
$ws.Range("E5").Value = 25
$ws.Range("F5").Value = "有"
$ws.Range("G5").Value = "一個月前"
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 4
$ws.Range("M5").Value = 4
$ws.Range("N5").Value = 5
$ws.Range("O5").Value = 3

# Blank respondent rows 6 through 15 (only numbered in column A for future entries)
$ws.Range("A6").Value = 3
$ws.Range("A7").Value = 4
$ws.Range("A8").Value = 5
$ws.Range("A9").Value = 6
$ws.Range("A10").Value = 7
$ws.Range("A11").Value = 8
$ws.Range("A12").Value = 9
$ws.Range("A13").Value = 10
$ws.Range("A14").Value = 11
$ws.Range("A15").Value = 12

# Summary/average row 16
$ws.Range("E16").Formula = "=AVERAGE(E4:E15)"
$ws.Range("G16").Value = "平均"
$ws.Range("H16:O16").Formula = "=AVERAGE(H4:H15)"

# Highlight the whole summary row with a bold red font on a light gold fill
$summaryRow = $ws.Range("A16:O16")
$summaryRow.Font.Bold = $true
$summaryRow.Font.Color = 255
$summaryRow.Interior.Color = 13431551

# Update print setup and the active selection to match the edited sheet
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.Range("E21").Select()
